$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the condition statement in A4: "status=1" -> "status=1 || status=2"
$ws.Range("A4").Value = "status=1 || status=2"

# Reflect the updated selected cell (cursor moved to A4 when last saved)
$ws.Range("A4").Select()
